$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stale "_GoBack" bookmark (left over from the previous save's
#    last-edit position). Word regenerates / drops this automatically when a
#    document is edited and re-saved; we drop it explicitly here so the
#    bookmarkStart/bookmarkEnd pair around the two recipe photos disappears.
# ---------------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}

# ---------------------------------------------------------------------------
# 2) Append new content at the very end of the document, after the final
#    "Enjoy your dessert!" paragraph:
#      - a blank paragraph
#      - a paragraph reading "This is the sentence I entered on 12th March 2020."
#        (with "th" as a superscript)
#      - a trailing blank paragraph
# ---------------------------------------------------------------------------

# Blank paragraph directly after "Enjoy your dessert!".
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.Text = "`r"

# Second paragraph break -- this paragraph will receive the sentence text.
# (Collapsing exactly at Content.End triggers an off-by-one in this host, so
# anchor one character before it instead.)
$endPos = $d.Content.End - 1
$r2 = $d.Range($endPos, $endPos)
$r2.InsertAfter("`r")

# Insert the whole sentence as one run (keeps paragraph-mark formatting
# such as snapToGrid/kern intact), then re-split "th" into a superscript run.
$endPos = $d.Content.End - 1
$r3 = $d.Range($endPos, $endPos)
$sentStart = $r3.Start
$sentence = "This is the sentence I entered on 12th March 2020."
$r3.InsertAfter($sentence)

$thRelStart = $sentence.IndexOf("12th") + 2
$thStart = $sentStart + $thRelStart
$thEnd = $thStart + 2
$rTh = $d.Range($thStart, $thEnd)
$rTh.Font.Superscript = $true

# Trailing blank paragraph.
$endPos = $d.Content.End
$r4 = $d.Range($endPos, $endPos)
$r4.Text = "`r"
